$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Status text: "Ready for handoff" -> "In Translation"
#    (Overview!E2, Overview!F2, zh-cn!C2, de-de!C2 all share this string)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# ---------------------------------------------------------------------
# 2) Narrow the "zh-cn" / "de-de" / "Status" columns.
#    ColumnWidth is expressed in characters (Calibri 11 / MDW 7), the
#    same unit Excel's own Format > Column Width dialog uses.
#    Overview: column E = zh-cn, column F = de-de
#    zh-cn / de-de sheets: column C = Status
# ---------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
